$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that flip from CHARTER -> STAY
$toStay = @(5, 6, 7, 9, 82, 85, 86, 93, 99, 170, 172)
foreach ($r in $toStay) {
    $ws.Range("A$r").Value = "STAY"
}

# Rows that flip from STAY -> CHARTER
$toCharter = @(95, 96, 100, 101, 109, 148)
foreach ($r in $toCharter) {
    $ws.Range("A$r").Value = "CHARTER"
}

# Update the active selection on the sheet view
$ws.Range("I14").Select()
